$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.652.68"
$ws.Range("E2").Value = "  +6.29%  "
$ws.Range("D3").Value = "2.639.81"
$ws.Range("E3").Value = "  +9.79%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'513.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.45%  "
$ws.Range("D6").Value = "'158.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "2.685.17"
$ws.Range("E9").Value = "  +10.93%  "
$ws.Range("D10").Value = "'6.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "'0.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.76%  "
$ws.Range("E12").Value = "  +4.21%  "
$ws.Range("D14").Value = "3.142.16"
$ws.Range("E14").Value = "  +11.13%  "
$ws.Range("D15").Value = "60.760.80"
$ws.Range("E15").Value = "  +6.59%  "
$ws.Range("D16").Value = "'21.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000141"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.22%  "
$ws.Range("D18").Value = "2.684.72"
$ws.Range("E18").Value = "  +10.81%  "
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").Value = "'349.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.79%  "
$ws.Range("D21").Value = "'10.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.48%  "
$ws.Range("D22").Value = "'6.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.26%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'60.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("D25").Value = "'0.425"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.91%  "
$ws.Range("D26").Value = "2.796.59"
$ws.Range("E26").Value = "  +10.85%  "
$ws.Range("E27").Value = "  +4.73%  "
$ws.Range("D28").Value = "'0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "0.0₃0878"
$ws.Range("E29").Value = "  +12.52%  "
$ws.Range("D30").Value = "'7.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.33%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'19.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.80%  "
$ws.Range("D33").Value = "'157.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.84%  "
$ws.Range("E34").Value = "  +3.88%  "
$ws.Range("E35").Value = "  +9.61%  "
$ws.Range("E36").Value = "  +10.28%  "
$ws.Range("D37").Value = "'1.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.65%  "
$ws.Range("D38").Value = "'313.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +16.65%  "
$ws.Range("D39").Value = "'1.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.15%  "
$ws.Range("D40").Value = "'0.863"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("D41").Value = "'0.846"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +32.40%  "
$ws.Range("D42").Value = "'3.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.28%  "
$ws.Range("D43").Value = "'35.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.50%  "
$ws.Range("D44").Value = "'0.647"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.85%  "
$ws.Range("D45").Value = "'0.0583"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.83%  "
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("E47").Value = "  +16.34%  "
$ws.Range("D48").Value = "'0.991"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'4.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.15%  "
$ws.Range("D50").Value = "2.078.51"
$ws.Range("E50").Value = "  +10.96%  "
$ws.Range("E51").Value = "  +3.36%  "
